$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "'17.56"
$ws.Range("C11").Value = "'5.26"
$ws.Range("D11").Value = "'22.82"

# Enterprises (% of total) - row 12
$ws.Range("B12").Value = "'74.89"
$ws.Range("C12").Value = "'22.43"
$ws.Range("D12").Value = "'97.31"

# Re-apply the original (default) style so the apostrophe-forced text
# entry above doesn't leave behind a "quote prefix" number-format style
# that wasn't part of the original formatting.
$ws.Range("B11:D12").Style = "Normal"
